$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 217, pushing existing rows 217-246 down to 218-247.
$ws.Rows.Item(217).Insert()

# Populate the newly inserted row 217 with the new weekly price record.
$ws.Cells.Item(217, 1).Value2 = 8
$ws.Cells.Item(217, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(217, 3).Value2 = "Coquimbo"
$ws.Cells.Item(217, 4).Value2 = 45127
$ws.Cells.Item(217, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(217, 5).Value2 = 4
$ws.Cells.Item(217, 6).Value2 = 100112040
$ws.Cells.Item(217, 7).Value2 = "Cilantro"
$ws.Cells.Item(217, 8).Value2 = "Sin especificar"
$ws.Cells.Item(217, 9).Value2 = "Primera"
$ws.Cells.Item(217, 10).Value2 = 2200
$ws.Cells.Item(217, 11).Value2 = 2500
$ws.Cells.Item(217, 12).Value2 = 3000
$ws.Cells.Item(217, 13).Value2 = 2750
$ws.Cells.Item(217, 14).Value2 = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(217, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(217, 16).Value2 = 1833
$ws.Cells.Item(217, 17).Value2 = 1.5
$ws.Cells.Item(217, 18).Value2 = "Hortaliza"
